# Auto-generated script applying cell value updates per the OOXML diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1290.72
$ws.Range("I15").Value = 1290.72
$ws.Range("K15").Value = 3872.16
$ws.Range("M15").Value = -3703.16
$ws.Range("H17").Value = 1039.9032
$ws.Range("J17").Value = 875.59015
$ws.Range("L17").Value = 2626.77045
$ws.Range("N17").Value = -2962.77045
$ws.Range("H47").Value = 9689
$ws.Range("I47").Value = 7033.5
$ws.Range("K47").Value = 7033.5
$ws.Range("M47").Value = -6061.5
$ws.Range("H62").Value = 2999.2
$ws.Range("I62").Value = 3000
$ws.Range("K62").Value = 3000
$ws.Range("M62").Value = -2376
$ws.Range("H65").Value = 2999.2
$ws.Range("I65").Value = 3000
$ws.Range("K65").Value = 15000
$ws.Range("M65").Value = -11880
$ws.Range("H121").Value = 150
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
$ws.Range("H131").Value = 3525.2727
$ws.Range("J131").Value = 4388
$ws.Range("L131").Value = 13164
$ws.Range("N131").Value = -23244
$ws.Range("H137").Value = 1665
$ws.Range("I137").Value = 1680.8334
$ws.Range("J137").Value = 1633.3334
$ws.Range("K137").Value = 5042.5002
$ws.Range("L137").Value = 4900.0002
$ws.Range("M137").Value = -2492.5002
$ws.Range("N137").Value = -10000.0002
$ws.Range("H138").Value = 3480.8057
$ws.Range("I138").Value = 3680.4119
$ws.Range("J138").Value = 3302.2104
$ws.Range("K138").Value = 11041.2357
$ws.Range("L138").Value = 9906.6312
$ws.Range("M138").Value = -5901.235700000001
$ws.Range("N138").Value = -20186.6312
$ws.Range("H141").Value = 1650578.9
$ws.Range("I141").Value = 2547796
$ws.Range("J141").Value = 5680.8335
$ws.Range("K141").Value = 7643388
$ws.Range("L141").Value = 17042.5005
$ws.Range("M141").Value = -7638208
$ws.Range("N141").Value = -27402.5005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5152.946
$ws.Range("I32").Value = 3815.9312
$ws.Range("J32").Value = 9999.625
$ws.Range("K32").Value = 3815.9312
$ws.Range("L32").Value = 9999.625
$ws.Range("M32").Value = -3528.9312
$ws.Range("N32").Value = -10573.625
$ws.Range("H61").Value = 4110
$ws.Range("I61").Value = 2770.9412
$ws.Range("J61").Value = 6955.5
$ws.Range("K61").Value = 2770.9412
$ws.Range("L61").Value = 6955.5
$ws.Range("M61").Value = -2558.9412
$ws.Range("N61").Value = -7379.5
$ws.Range("H74").Value = 1298.8292
$ws.Range("I74").Value = 855.4865
$ws.Range("J74").Value = 5399.75
$ws.Range("K74").Value = 855.4865
$ws.Range("L74").Value = 5399.75
$ws.Range("M74").Value = 18.51350000000002
$ws.Range("N74").Value = -7147.75
$ws.Range("H77").Value = 1298.8292
$ws.Range("I77").Value = 855.4865
$ws.Range("J77").Value = 5399.75
$ws.Range("K77").Value = 4277.4325
$ws.Range("L77").Value = 26998.75
$ws.Range("M77").Value = 90.56750000000011
$ws.Range("N77").Value = -35734.75
$ws.Range("H110").Value = 3276.7144
$ws.Range("I110").Value = 1852.75
$ws.Range("J110").Value = 5175.3335
$ws.Range("K110").Value = 1852.75
$ws.Range("L110").Value = 5175.3335
$ws.Range("M110").Value = 192.25
$ws.Range("N110").Value = -9265.333500000001
$ws.Range("H132").Value = 1764.125
$ws.Range("I132").Value = 1073.4286
$ws.Range("K132").Value = 3220.2858
$ws.Range("M132").Value = -690.2857999999997
$ws.Range("H136").Value = 4110
$ws.Range("I136").Value = 2770.9412
$ws.Range("J136").Value = 6955.5
$ws.Range("K136").Value = 8312.8236
$ws.Range("L136").Value = 20866.5
$ws.Range("M136").Value = -5762.8236
$ws.Range("N136").Value = -25966.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H107").Value = 3044.2856
$ws.Range("I107").Value = 3044.2856
$ws.Range("K107").Value = 3044.2856
$ws.Range("M107").Value = -1124.2856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 537
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
$ws.Range("H31").Value = 3782.6428
$ws.Range("I31").Value = 2845.5715
$ws.Range("J31").Value = 4719.7144
$ws.Range("K31").Value = 2845.5715
$ws.Range("L31").Value = 4719.7144
$ws.Range("M31").Value = -2550.5715
$ws.Range("N31").Value = -5309.7144
$ws.Range("H34").Value = 3782.6428
$ws.Range("I34").Value = 2845.5715
$ws.Range("J34").Value = 4719.7144
$ws.Range("K34").Value = 2845.5715
$ws.Range("L34").Value = 4719.7144
$ws.Range("M34").Value = -2643.5715
$ws.Range("N34").Value = -5123.7144
$ws.Range("H99").Value = 2652.3333
$ws.Range("J99").Value = 3228.5
$ws.Range("L99").Value = 3228.5
$ws.Range("N99").Value = -6224.5
$ws.Range("H126").Value = 2652.3333
$ws.Range("J126").Value = 3228.5
$ws.Range("L126").Value = 9685.5
$ws.Range("N126").Value = -14625.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 173
$ws.Range("I4").Value = 173
$ws.Range("K4").Value = 519
$ws.Range("M4").Value = -407
$ws.Range("H129").Value = 33913.363
$ws.Range("J129").Value = 41294.945
$ws.Range("L129").Value = 123884.835
$ws.Range("N129").Value = -133884.835
$ws.Range("H131").Value = 754.69385
$ws.Range("I131").Value = 447.1111
$ws.Range("J131").Value = 785.7977
$ws.Range("K131").Value = 1341.3333
$ws.Range("L131").Value = 2357.3931
$ws.Range("M131").Value = 3698.6667
$ws.Range("N131").Value = -12437.3931
$ws.Range("H140").Value = 1721.2307
$ws.Range("I140").Value = 977.8461
$ws.Range("J140").Value = 2464.6155
$ws.Range("K140").Value = 2933.5383
$ws.Range("L140").Value = 7393.8465
$ws.Range("M140").Value = 2246.4617
$ws.Range("N140").Value = -17753.8465

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("H10").Value = 1667333
$ws.Range("I10").Value = 2500499.5
$ws.Range("K10").Value = 2500499.5
$ws.Range("M10").Value = -2500330.5
$ws.Range("H80").Value = 2974
$ws.Range("I80").Value = 2873.4285
$ws.Range("K80").Value = 2873.4285
$ws.Range("M80").Value = -1875.4285
$ws.Range("H83").Value = 2974
$ws.Range("I83").Value = 2873.4285
$ws.Range("K83").Value = 14367.1425
$ws.Range("M83").Value = -9375.1425
$ws.Range("H102").Value = 4493.1665
$ws.Range("I102").Value = 4791.6
$ws.Range("K102").Value = 4791.6
$ws.Range("M102").Value = -3169.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2433.111
$ws.Range("I22").Value = 3666.6667
$ws.Range("J22").Value = 1816.3334
$ws.Range("K22").Value = 3666.6667
$ws.Range("L22").Value = 1816.3334
$ws.Range("M22").Value = -3371.6667
$ws.Range("N22").Value = -2406.3334
$ws.Range("H27").Value = 2433.111
$ws.Range("I27").Value = 3666.6667
$ws.Range("J27").Value = 1816.3334
$ws.Range("K27").Value = 3666.6667
$ws.Range("L27").Value = 1816.3334
$ws.Range("M27").Value = -3559.6667
$ws.Range("N27").Value = -2030.3334
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H132").Value = 2118.0476
$ws.Range("I132").Value = 1912.6
$ws.Range("K132").Value = 5737.799999999999
$ws.Range("M132").Value = -3207.799999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 80005
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()

Write-Output "Updated cells: sets=194 clears=7"